# Fixing errors in example upload files.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Service Contacts sheet: add a custom width for column A and move the
# active selection to D3.
# ---------------------------------------------------------------------------
$wsSvc = $wb.Worksheets.Item("Service Contacts")
$wsSvc.Columns.Item(1).ColumnWidth = 13.6666667
$wsSvc.Range("D3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Practitioners sheet: widen a few columns, add a new practitioner record
# (row 6) and update the active selection.
# ---------------------------------------------------------------------------
$wsPrac = $wb.Worksheets.Item("Practitioners")

$wsPrac.Columns.Item(1).ColumnWidth = 13.8333333
$wsPrac.Columns.Item(3).ColumnWidth = 12.1666667
$wsPrac.Columns.Item(6).ColumnWidth = 12.0

$wsPrac.Cells.Item(6, 1).Value = "PHN999:NFP02"
$wsPrac.Cells.Item(6, 2).Value = "P01"
$wsPrac.Cells.Item(6, 3).Value = 8
$wsPrac.Cells.Item(6, 4).Value = 1
$wsPrac.Cells.Item(6, 5).Value = 1973
$wsPrac.Cells.Item(6, 6).Value = 2
$wsPrac.Cells.Item(6, 7).Value = 1
$wsPrac.Cells.Item(6, 8).Value = 1
$wsPrac.Cells.Item(6, 9).Value = "tag1"

$wsPrac.Range("G1:G1048576").Select() | Out-Null

# Restore the originally active sheet/tab (SDQ) so the workbook-level
# "activeTab" / per-sheet "tabSelected" state is left unchanged, matching
# only the per-sheet selection updates made above.
$wb.Worksheets.Item("SDQ").Activate() | Out-Null
